$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text (matches
# the source workbook, which keeps these as literal strings) while leaving
# the cell style index untouched (format -> Text, set value, format back to
# General, then restore the Normal cell style so no "s" attribute lingers).
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "29.873.96"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.891.78"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "0.7760"
$ws.Range("E5").Value = "  -1.67%  "
Set-TextValue "D6" "244.28"
$ws.Range("E6").Value = "  +0.90%  "
Set-TextValue "D7" "1.000"
$ws.Range("E7").Value = "  -0.15%  "
Set-TextValue "D8" "0.3141"
$ws.Range("E8").Value = "  -1.44%  "
Set-TextValue "D9" "0.07395"
$ws.Range("E9").Value = "  +4.64%  "
$ws.Range("E10").Value = "  -1.92%  "
Set-TextValue "D11" "0.08149"
$ws.Range("E11").Value = "  +1.17%  "
Set-TextValue "D12" "0.7664"
$ws.Range("E12").Value = "  -0.64%  "
Set-TextValue "D13" "5.485"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").Value = "1.868.64"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("E15").Value = "  +0.29%  "
Set-TextValue "D16" "6.216"
$ws.Range("E16").Value = "  +5.21%  "
$ws.Range("D17").Value = "29.862.65"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +0.77%  "
Set-TextValue "D19" "245.28"
$ws.Range("E19").Value = "  +0.67%  "
Set-TextValue "D20" "0.000007864"
$ws.Range("E20").Value = "  +2.00%  "
Set-TextValue "D21" "0.9998"
$ws.Range("E21").Value = "  -0.15%  "
Set-TextValue "D22" "8.126"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "2.111.49"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  -0.20%  "
Set-TextValue "D25" "0.1567"
$ws.Range("E25").Value = "  -3.33%  "
Set-TextValue "D26" "9.435"
$ws.Range("E26").Value = "  +1.41%  "
Set-TextValue "D27" "162.21"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  +5.49%  "
Set-TextValue "D31" "1.546"
$ws.Range("E31").Value = "  +0.76%  "
Set-TextValue "D32" "4.499"
$ws.Range("E32").Value = "  +2.20%  "
Set-TextValue "D33" "0.05603"
$ws.Range("E33").Value = "  -0.40%  "
Set-TextValue "D34" "4.100"
$ws.Range("E34").Value = "  -0.06%  "
Set-TextValue "D35" "1.252"
$ws.Range("E35").Value = "  -1.30%  "
Set-TextValue "D36" "0.7587"
$ws.Range("E36").Value = "  +3.04%  "
Set-TextValue "D37" "0.9984"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue "D38" "2.649"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  +0.41%  "
Set-TextValue "D40" "2.788"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "1.152.21"
$ws.Range("E41").Value = "  +12.93%  "
$ws.Range("E42").Value = "  +2.99%  "
Set-TextValue "D43" "0.4464"
$ws.Range("E43").Value = "  +0.40%  "
Set-TextValue "D44" "5.968"
Set-TextValue "D45" "0.8537"
$ws.Range("E45").Value = "  +1.00%  "
Set-TextValue "D46" "1.905"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("E47").Value = "  -0.09%  "
Set-TextValue "D48" "3.122"
$ws.Range("E48").Value = "  +5.60%  "
Set-TextValue "D49" "101.92"
$ws.Range("E49").Value = "  -0.48%  "
Set-TextValue "D50" "9.869"
$ws.Range("E50").Value = "  -0.80%  "
Set-TextValue "D51" "7.514"
$ws.Range("E51").Value = "  +0.49%  "
